# Update "Förändrad" (Changed) date column C for rows 2-42 from 45711 to 45712
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C42").Value = 45712
